$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete obsolete rows 8-10 (data now only spans rows 2-7)
$ws.Range("A8:A10").EntireRow.Delete()

# Overwrite rows 2-7 with the updated TPM-derived values
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icam4"
$ws.Range("C2").Value = "Itgb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3732763333333333
$ws.Range("H2").Value = 1.119829
$ws.Range("I2").Value = 0.3554258969843855
$ws.Range("J2").Value = 0.3554258969843855
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.02097464602088889
$ws.Range("R2").Value = 0.188771814188
$ws.Range("S2").Value = 0.1397316937543334
$ws.Range("T2").Value = 0.1397316937543334

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icam4"
$ws.Range("C3").Value = "Itgb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3732763333333333
$ws.Range("H3").Value = 1.119829
$ws.Range("I3").Value = 0.3554258969843855
$ws.Range("J3").Value = 0.3554258969843855
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08673766666666667
$ws.Range("N3").Value = 0.260213
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("Q3").Value = 0.03237711817522222
$ws.Range("R3").Value = 0.291394063577
$ws.Range("S3").Value = 0.2156942032300521
$ws.Range("T3").Value = 0.2156942032300521

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Icam4"
$ws.Range("C4").Value = "Itgb2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.668317
$ws.Range("H4").Value = 2.004951
$ws.Range("I4").Value = 0.6363574327729865
$ws.Range("J4").Value = 0.6363574327729865
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05619066666666667
$ws.Range("N4").Value = 0.168572
$ws.Range("O4").Value = 0.3931387525216601
$ws.Range("P4").Value = 0.39313875252166
$ws.Range("Q4").Value = 0.03755317777466667
$ws.Range("R4").Value = 0.337978599972
$ws.Range("S4").Value = 0.2501767672782581
$ws.Range("T4").Value = 0.250176767278258

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icam4"
$ws.Range("C5").Value = "Itgb2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.668317
$ws.Range("H5").Value = 2.004951
$ws.Range("I5").Value = 0.6363574327729865
$ws.Range("J5").Value = 0.6363574327729865
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08673766666666667
$ws.Range("N5").Value = 0.260213
$ws.Range("O5").Value = 0.60686124747834
$ws.Range("P5").Value = 0.60686124747834
$ws.Range("Q5").Value = 0.05796825717366667
$ws.Range("R5").Value = 0.5217143145630001
$ws.Range("S5").Value = 0.3861806654947285
$ws.Range("T5").Value = 0.3861806654947285

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Icam4"
$ws.Range("C6").Value = "Itgb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.008629333333333334
$ws.Range("H6").Value = 0.025888
$ws.Range("I6").Value = 0.008216670242627913
$ws.Range("J6").Value = 0.008216670242627911
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05619066666666667
$ws.Range("N6").Value = 0.168572
$ws.Range("O6").Value = 0.3931387525216601
$ws.Range("P6").Value = 0.39313875252166
$ws.Range("Q6").Value = 0.0004848879928888889
$ws.Range("R6").Value = 0.004363991936
$ws.Range("S6").Value = 0.003230291489068584
$ws.Range("T6").Value = 0.003230291489068583

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Icam4"
$ws.Range("C7").Value = "Itgb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.008629333333333334
$ws.Range("H7").Value = 0.025888
$ws.Range("I7").Value = 0.008216670242627913
$ws.Range("J7").Value = 0.008216670242627911
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08673766666666667
$ws.Range("N7").Value = 0.260213
$ws.Range("O7").Value = 0.60686124747834
$ws.Range("P7").Value = 0.60686124747834
$ws.Range("Q7").Value = 0.0007484882382222224
$ws.Range("R7").Value = 0.006736394144000001
$ws.Range("S7").Value = 0.00498637875355933
$ws.Range("T7").Value = 0.004986378753559329
